$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format / style) from the last existing data
# row (row 16) down onto the two new rows being appended.
$ws.Range("A16").Copy()
$ws.Range("A17:A18").PasteSpecial(-4122)  # xlPasteFormats

# Row 17
$ws.Range("A17").Value = 45748
$ws.Range("B17").Value = 0.23694
$ws.Range("C17").Value = 0.21437
$ws.Range("D17").Value = 0.34707
$ws.Range("E17").Value = 0.20162
$ws.Range("F17").Value = 0.04546

# Row 18
$ws.Range("A18").Value = 45778
$ws.Range("B18").Value = 0.25321
$ws.Range("C18").Value = 0.21259
$ws.Range("D18").Value = 0.36264
$ws.Range("E18").Value = 0.17156
$ws.Range("F18").Value = 0.04546

$wb.Save()
